# This script updates the "想去人数" (want-to-go count) column F values
# on sheets 1 (展览), 2 (演出), and 4 (全部类型) to reflect refreshed
# data as published at commit 456a3b4 (gh-pages data regeneration).
# Sheet 3 (本地生活) has no changes in this update.

$wb = $excel.ActiveWorkbook

# --- Worksheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7853
$ws.Range("F3").Value = 7998
$ws.Range("F5").Value = 46
$ws.Range("F6").Value = 6997
$ws.Range("F7").Value = 3461
$ws.Range("F9").Value = 3774
$ws.Range("F12").Value = 61
$ws.Range("F14").Value = 106
$ws.Range("F15").Value = 491
$ws.Range("F17").Value = 95
$ws.Range("F18").Value = 344
$ws.Range("F21").Value = 345
$ws.Range("F26").Value = 515
$ws.Range("F27").Value = 1584
$ws.Range("F30").Value = 2902
$ws.Range("F31").Value = 2079
$ws.Range("F33").Value = 69
$ws.Range("F36").Value = 4010
$ws.Range("F37").Value = 400
$ws.Range("F40").Value = 936
$ws.Range("F41").Value = 701
$ws.Range("F42").Value = 113
$ws.Range("F44").Value = 1537
$ws.Range("F45").Value = 256

# --- Worksheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 431
$ws.Range("F5").Value = 56
$ws.Range("F7").Value = 105
$ws.Range("F11").Value = 104
$ws.Range("F14").Value = 19
$ws.Range("F15").Value = 511

# --- Worksheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 7853
$ws.Range("F5").Value = 7998
$ws.Range("F7").Value = 46
$ws.Range("F8").Value = 6997
$ws.Range("F9").Value = 3461
$ws.Range("F11").Value = 3774
$ws.Range("F14").Value = 61
$ws.Range("F15").Value = 106
$ws.Range("F17").Value = 56
$ws.Range("F18").Value = 95
$ws.Range("F21").Value = 105
$ws.Range("F23").Value = 345
$ws.Range("F28").Value = 515
$ws.Range("F29").Value = 1584
$ws.Range("F32").Value = 2902
$ws.Range("F33").Value = 2079
$ws.Range("F35").Value = 69
$ws.Range("F37").Value = 104
$ws.Range("F38").Value = 4010
$ws.Range("F39").Value = 400
$ws.Range("F41").Value = 19
$ws.Range("F43").Value = 936
$ws.Range("F44").Value = 701
$ws.Range("F45").Value = 1537
$ws.Range("F46").Value = 256

